# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Fri Jul 14 11:33:43 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while keeping it stored as TEXT (matches the
# original inline-string cells) even when the text looks like a plain number,
# e.g. "0.9978" or "165.59". Without this, Excel would silently convert such
# strings to real numbers. ClearFormats() afterwards drops the temporary '@'
# text number-format again so the cell keeps its original (default) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range('D2').Value = '31.153.34'
$ws.Range('E2').Value = '  +1.87%  '
# Row 3
$ws.Range('D3').Value = '1.988.05'
$ws.Range('E3').Value = '  +5.51%  '
# Row 4
Set-TextValue $ws.Range('D4') '0.9978'
$ws.Range('E4').Value = '  -0.24%  '
# Row 5
Set-TextValue $ws.Range('D5') '0.8006'
$ws.Range('E5').Value = '  +69.00%  '
# Row 6
Set-TextValue $ws.Range('D6') '254.87'
$ws.Range('E6').Value = '  +3.38%  '
# Row 7
Set-TextValue $ws.Range('D7') '0.9976'
$ws.Range('E7').Value = '  -0.25%  '
# Row 8
Set-TextValue $ws.Range('D8') '0.3499'
# Row 9
Set-TextValue $ws.Range('D9') '28.15'
$ws.Range('E9').Value = '  +25.98%  '
# Row 10
Set-TextValue $ws.Range('D10') '0.06986'
$ws.Range('E10').Value = '  +6.79%  '
# Row 11
Set-TextValue $ws.Range('D11') '0.8457'
$ws.Range('E11').Value = '  +8.84%  '
# Row 12
Set-TextValue $ws.Range('D12') '0.08176'
$ws.Range('E12').Value = '  +4.81%  '
# Row 13
$ws.Range('D13').Value = '1.990.43'
$ws.Range('E13').Value = '  +5.65%  '
# Row 14
Set-TextValue $ws.Range('D14') '100.43'
$ws.Range('E14').Value = '  -0.47%  '
# Row 15
Set-TextValue $ws.Range('D15') '5.619'
$ws.Range('E15').Value = '  +6.81%  '
# Row 16
Set-TextValue $ws.Range('D16') '15.43'
$ws.Range('E16').Value = '  +16.69%  '
# Row 17
Set-TextValue $ws.Range('D17') '273.12'
$ws.Range('E17').Value = '  -4.28%  '
# Row 18
$ws.Range('D18').Value = '31.152.22'
$ws.Range('E18').Value = '  +1.94%  '
# Row 19
$ws.Range('E19').Value = '  +9.80%  '
# Row 20
Set-TextValue $ws.Range('D20') '0.000007930'
$ws.Range('E20').Value = '  +5.31%  '
# Row 21
$ws.Range('D21').Value = '2.252.13'
$ws.Range('E21').Value = '  +5.94%  '
# Row 22
Set-TextValue $ws.Range('D22') '0.9984'
$ws.Range('E22').Value = '  -0.18%  '
# Row 23
Set-TextValue $ws.Range('D23') '0.9971'
$ws.Range('E23').Value = '  -0.30%  '
# Row 24
Set-TextValue $ws.Range('D24') '7.050'
$ws.Range('E24').Value = '  +9.53%  '
# Row 25
Set-TextValue $ws.Range('D25') '9.989'
$ws.Range('E25').Value = '  +8.93%  '
# Row 26
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D26') '0.1514'
$ws.Range('E26').Value = '  +56.00%  '
# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D27') '165.59'
$ws.Range('E27').Value = '  +2.01%  '
# Row 28
Set-TextValue $ws.Range('D28') '19.90'
$ws.Range('E28').Value = '  +4.11%  '
# Row 29
Set-TextValue $ws.Range('D29') '2.344'
$ws.Range('E29').Value = '  +22.49%  '
# Row 30
Set-TextValue $ws.Range('D30') '1.596'
$ws.Range('E30').Value = '  +6.12%  '
# Row 31
Set-TextValue $ws.Range('D31') '1.355'
$ws.Range('E31').Value = '  +1.70%  '
# Row 32
Set-TextValue $ws.Range('D32') '4.579'
$ws.Range('E32').Value = '  +7.47%  '
# Row 33
Set-TextValue $ws.Range('D33') '4.409'
$ws.Range('E33').Value = '  +5.19%  '
# Row 34
Set-TextValue $ws.Range('D34') '0.05245'
$ws.Range('E34').Value = '  +8.21%  '
# Row 35
Set-TextValue $ws.Range('D35') '0.7796'
$ws.Range('E35').Value = '  +11.79%  '
# Row 36
Set-TextValue $ws.Range('D36') '1.219'
$ws.Range('E36').Value = '  +7.82%  '
# Row 37
Set-TextValue $ws.Range('D37') '2.763'
$ws.Range('E37').Value = '  +0.24%  '
# Row 38
Set-TextValue $ws.Range('D38') '0.9975'
$ws.Range('E38').Value = '  -0.22%  '
# Row 40
Set-TextValue $ws.Range('D40') '2.881'
$ws.Range('E40').Value = '  -0.61%  '
# Row 41
Set-TextValue $ws.Range('D41') '6.653'
$ws.Range('E41').Value = '  +5.80%  '
# Row 42
Set-TextValue $ws.Range('D42') '79.75'
$ws.Range('E42').Value = '  +4.96%  '
# Row 43
Set-TextValue $ws.Range('D43') '0.4668'
$ws.Range('E43').Value = '  +9.65%  '
# Row 44
Set-TextValue $ws.Range('D44') '2.127'
$ws.Range('E44').Value = '  +7.29%  '
# Row 45
Set-TextValue $ws.Range('D45') '0.8527'
$ws.Range('E45').Value = '  +2.75%  '
# Row 46
Set-TextValue $ws.Range('D46') '104.70'
$ws.Range('E46').Value = '  +3.06%  '
# Row 47
Set-TextValue $ws.Range('D47') '0.9975'
$ws.Range('E47').Value = '  -0.23%  '
# Row 48
Set-TextValue $ws.Range('D48') '7.673'
$ws.Range('E48').Value = '  +9.22%  '
# Row 49
Set-TextValue $ws.Range('D49') '9.864'
$ws.Range('E49').Value = '  -0.44%  '
# Row 50
$ws.Range('E50').Value = '  +4.50%  '
# Row 51
$ws.Range('E51').Value = '  +8.46%  '
